$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(13,2).Font.Name = "Arial"
$ws.Cells.Item(13,2).Font.Name = "Times New Roman"
